# Robo-advisor config: tweak the PortfolioIndex equity glidepath values
# (each step down from 0.08/0.18/.../0.98 to 0.05/0.15/.../0.95), and leave
# the workbook focused on the PortfolioIndex sheet/selection as the author
# last had it.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("PortfolioIndex")

$ws2.Range("B2").Value = 0.05
$ws2.Range("B3").Value = 0.15
$ws2.Range("B4").Value = 0.25
$ws2.Range("B5").Value = 0.35
$ws2.Range("B6").Value = 0.45
$ws2.Range("B7").Value = 0.55
$ws2.Range("B8").Value = 0.65
$ws2.Range("B9").Value = 0.75
$ws2.Range("B10").Value = 0.85
$ws2.Range("B11").Value = 0.95

# Make PortfolioIndex the active/selected tab and move its selection to L9
# (this also clears tabSelected from the Glidepath sheet, which was the
# previously active tab).
$ws2.Activate()
$ws2.Range("L9").Select()
